# Updates cryptos list values (Price and Volume(1h) columns) per commit
# "Updated cryptos list on Sun May 21 05:26:43 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.368.29"
$ws.Range("E2").Value = "  +1.70%  "
$ws.Range("D3").Value = "1.834.16"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.90%  "
$ws.Range("D5").Value = "'314.19"
$ws.Range("E5").Value = "  +1.72%  "
$ws.Range("D7").Value = "'0.4747"
$ws.Range("E7").Value = "  +2.15%  "
$ws.Range("D8").Value = "'0.3695"
$ws.Range("E8").Value = "  +1.04%  "
$ws.Range("D9").Value = "'0.07463"
$ws.Range("E9").Value = "  +1.36%  "
$ws.Range("D10").Value = "'0.8870"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D12").Value = "1.882.84"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("D13").Value = "'0.07348"
$ws.Range("E13").Value = "  +3.61%  "
$ws.Range("E14").Value = "  +1.72%  "
$ws.Range("D15").Value = "'93.09"
$ws.Range("E15").Value = "  +1.91%  "
$ws.Range("D16").Value = "'6.588"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").Value = "'0.000008826"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.90%  "
$ws.Range("D20").Value = "27.687.68"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").Value = "'14.80"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").Value = "'5.317"
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +0.79%  "
$ws.Range("D24").Value = "2.111.11"
$ws.Range("E24").Value = "  +4.47%  "
$ws.Range("D25").Value = "'1.892"
$ws.Range("E25").Value = "  -0.16%  "
$ws.Range("D26").Value = "'151.92"
$ws.Range("E26").Value = "  +0.69%  "
$ws.Range("D27").Value = "'18.64"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").Value = "'5.244"
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("D30").Value = "'117.63"
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("D31").Value = "'0.09006"
$ws.Range("E31").Value = "  +1.40%  "
$ws.Range("D32").Value = "'0.7577"
$ws.Range("E32").Value = "  +0.30%  "
$ws.Range("D33").Value = "'1.177"
$ws.Range("E33").Value = "  +1.69%  "
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").Value = "'2.947"
$ws.Range("E35").Value = "  +1.43%  "
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").Value = "'1.104"
$ws.Range("E37").Value = "  +1.70%  "
$ws.Range("D38").Value = "'0.05355"
$ws.Range("E38").Value = "  +1.29%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'2.987"
$ws.Range("E40").Value = "  +0.11%  "
$ws.Range("D41").Value = "'7.328"
$ws.Range("E41").Value = "  +1.07%  "
$ws.Range("D42").Value = "'2.402"
$ws.Range("E42").Value = "  +3.27%  "
$ws.Range("D43").Value = "'0.5325"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'0.1660"
$ws.Range("E44").Value = "  +0.50%  "
$ws.Range("D45").Value = "'8.505"
$ws.Range("E45").Value = "  +1.09%  "
$ws.Range("D46").Value = "'0.4922"
$ws.Range("E46").Value = "  +1.09%  "
$ws.Range("D47").Value = "'10.55"
$ws.Range("E47").Value = "  +1.85%  "
$ws.Range("E48").Value = "  +1.02%  "
$ws.Range("D49").Value = "'105.10"
$ws.Range("E49").Value = "  +1.81%  "
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Value = "'0.06295"
$ws.Range("E51").Value = "  -0.02%  "
